$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): F2 5136 -> 5155, F4 902 -> 904
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5155
$wsExhibit.Range("F4").Value = 904

# Sheet "演出" (shows): F2 1 -> 2
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 2

# Sheet "全部类型" (all types): F2 5136 -> 5155, F4 902 -> 904, F5 1 -> 2
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5155
$wsAll.Range("F4").Value = 904
$wsAll.Range("F5").Value = 2
